$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the test data for the SISE "generación de siniestro" validation.
# The claim/policy number and the claim date used for the second scenario are
# being refreshed; both are plain text values (leading apostrophe keeps the
# numeric-looking policy number and the dd/mm/yyyy date stored as text,
# matching the original quotePrefix-text style already on these cells).
$ws.Range("G2").Value = "'26/04/2021"
$ws.Range("E2").Value = "'11111003159"

# Move/save the active selection to match where the user left the cursor.
$ws.Range("E9").Select()
